$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D17").Value = 2
$ws.Range("E17").Value = "2026-02-19T08:52:56.717571+00:00"
$ws.Range("H17").Value = 7
$ws.Range("L17").Value = "[191121, 191156, 191185, 191125, 191187, 191205, 191213]"
